$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "lstm_5"
$ws.Range("E10").Value = "(None, 7, 1)"
$ws.Range("B31").Value = "lstm_6"
$ws.Range("B58").Value = "dense_3"
